$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Edit 1: Insert a new list item "(5) ..." right after the "(4) Press
# Enter" step, describing the install-confirmation messages users will
# see in their command prompt.
# --------------------------------------------------------------------
$pressEnter = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "(4) Press Enter" + [char]13) {
        $pressEnter = $p
    }
}

if ($pressEnter -ne $null) {
    $pressEnter.Range.InsertParagraphAfter()
    $newPara = $pressEnter.Next()
    $newPara.Range.Text = "(5) Your command prompt will issue a series of message indicating whether each of the required packages has successfully installed. If you a receive a failure message for a specific package, try installing that package separately."
}

# --------------------------------------------------------------------
# Edit 2: Fix a doubled space in "The next  screen will help you ..."
# --------------------------------------------------------------------
$d.Content.Find.Execute("The next  screen will help you", $true, $false, $false, $false, $false, $true, 1, $false, "The next screen will help you", 2) | Out-Null
